$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Content updates
# ---------------------------------------------------------------------------
# Existing risk #2 ("Arduíno queimar") gets a new "Como?" action.
$ws.Range("G3").Value = "Todos os membros verificarem a conexão antes de ligar o arduíno"

# Existing risk #4 ("Upar arquivo errado no Git Hub") gets a new "Como?" action.
$ws.Range("G5").Value = "Manter arquivo dentro da pasta do projeto com instruções de como mexer no git"

# Existing risk #5 ("Modelar o BD de maneira incorreta") gets a new "Como?" action.
$ws.Range("G6").Value = "Manter contato constante com os professores, buscando feedback"

# Existing risk #3 ("Problema de conexão do BD com o site") action text is reworded.
$ws.Range("G4").Value = "Caso haja erro com a conexão, repetiremos os processos até solucionar o erro."

# New row 7: risk #6 "Botões e links quebrados no site".
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Botões e links quebrados no site"
$ws.Range("C7").Value = "Provável"
$ws.Range("D7").Value = "Alto"
$ws.Range("F7").Value = "Evitar"
$ws.Range("G7").Value = "Testar todas as funcionalidades do site antes da apresentação"

# E/H columns exist (empty, styled) on every data row, including the new row 7.
$ws.Range("E2:E7").ClearContents()
$ws.Range("H2:H7").ClearContents()

# ---------------------------------------------------------------------------
# 2) Formatting: start from a clean slate, then rebuild consistently.
# ---------------------------------------------------------------------------
$ws.Cells.ClearFormats()

$used = $ws.Range("A1:H7")
$used.Font.Name = "Arial"
$used.Font.Size = 12
$used.WrapText = $true
$used.VerticalAlignment = -4108

# Header row + the Probabilidade/Impacto columns are centered.
$centered = $ws.Range("A1:H1,C2:D7")
$centered.HorizontalAlignment = -4108

# Everything else in the body is left-aligned.
$ws.Range("A2:B7").HorizontalAlignment = -4131
$ws.Range("E2:H7").HorizontalAlignment = -4131

# The "Manter arquivo..." action (G5) is underlined.
$ws.Range("G5").Font.Underline = 2

# ---------------------------------------------------------------------------
# 3) Row heights / column widths
# ---------------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 38
$ws.Rows.Item(4).RowHeight = 51
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 34
$ws.Rows.Item(7).RowHeight = 34

$ws.Columns.Item(1).ColumnWidth = 7.166666666666667
$ws.Columns.Item(2).ColumnWidth = 36.666666666666664
$ws.Columns.Item(3).ColumnWidth = 14.998697916666666
$ws.Columns.Item(4).ColumnWidth = 13.166666666666666
$ws.Columns.Item(5).ColumnWidth = 11.330729166666666
$ws.Columns.Item(6).ColumnWidth = 9.666666666666666
$ws.Columns.Item(7).ColumnWidth = 36.498697916666664
$ws.Columns.Item(8).ColumnWidth = 22.830729166666668

# ---------------------------------------------------------------------------
# 4) Selection
# ---------------------------------------------------------------------------
$ws.Range("H6").Select()

Write-Output "done"
